$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B through AD hold the match-record data (column A is just a
# positional index and must stay untouched). For each pair of rows below,
# the two records were swapped (same match ids/odds, rows exchanged).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$rowPairs = @(
    @(110, 111),
    @(224, 225),
    @(231, 232),
    @(237, 238),
    @(249, 250)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = $col + $r1
        $addr2 = $col + $r2

        $range1 = $ws.Range($addr1)
        $range2 = $ws.Range($addr2)

        $val1 = $range1.Value2
        $val2 = $range2.Value2

        $range1.Value = $val2
        $range2.Value = $val1
    }
}
